# Generate Report for Handback
# The 996f345d-... file has now been handed back (in sync with en-US),
# so it moves to the top row of each sheet's table and gains handback
# file/datetime information. The 089084b7-... file remains "Ready for
# handoff" and moves to the second row.

$wb = $excel.ActiveWorkbook

# ---- URLs reused from the existing hyperlink relationships ----
$mdUrl_996 = "https://github.com/OpenLocalizationTest/oltest/blob/3b2e8129a28b74df82b18589fa076b11462cea47/e2e/996f345d-76ed-4a36-aca8-7897bfd7bef4.md"
$mdUrl_089 = "https://github.com/OpenLocalizationTest/oltest/blob/85478c745bd0d5db2800773d6f23f9abd8345986/e2e/089084b7-a469-4f48-856b-98300d6a7fc4.md"

$xlfUrl_996_zhcn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74d458d2da15eab048fd23eec18b3cb3abf9c817/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.zh-cn.xlf"
$xlfUrl_089_zhcn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/29fb1cbc2de41897035ee53ffaf45721a48358e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/089084b7-a469-4f48-856b-98300d6a7fc4.8ddaf840d9f6259ca43465a674b6b03e27c7b072.zh-cn.xlf"

$xlfUrl_996_dede = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/492b345300e02bf51d1fd7cab08091c74318754c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.de-de.xlf"
$xlfUrl_089_dede = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1eef84c32f57a4118fe30615b9e9781666cdfb24/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/089084b7-a469-4f48-856b-98300d6a7fc4.8ddaf840d9f6259ca43465a674b6b03e27c7b072.de-de.xlf"

$md996 = "996f345d-76ed-4a36-aca8-7897bfd7bef4.md"
$md089 = "089084b7-a469-4f48-856b-98300d6a7fc4.md"

# ==================== Overview sheet ====================
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Cells.Item(2,1).Value = $md996
$ws.Cells.Item(2,2).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,4).Value = "2016-30-18 16:30:55"

$ws.Cells.Item(3,1).Value = $md089
$ws.Cells.Item(3,2).Value = "Ready for handoff"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "2016-30-18 16:30:26"

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_996, "", "", $md996) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_089, "", "", $md089) | Out-Null

# ==================== zh-cn sheet ====================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Cells.Item(2,1).Value = $md996
$ws.Cells.Item(2,2).Value = ".md"
$ws.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,4).Value = "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.zh-cn.xlf"
$ws.Cells.Item(2,5).Value = "2016-03-18 16:30:52"
$ws.Cells.Item(2,6).Value = $md996
$ws.Cells.Item(2,7).Value = "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.zh-cn.xlf"
$ws.Cells.Item(2,8).Value = "2016-03-18 16:31:39"
$ws.Cells.Item(2,9).Value = "Include"

$ws.Cells.Item(3,1).Value = $md089
$ws.Cells.Item(3,2).Value = ".md"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "089084b7-a469-4f48-856b-98300d6a7fc4.8ddaf840d9f6259ca43465a674b6b03e27c7b072.zh-cn.xlf"
$ws.Cells.Item(3,5).Value = "2016-03-18 16:30:23"
$ws.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(3,9).Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_996, "", "", $md996) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl_996, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrl_996_zhcn, "", "", "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl_996, "", "", $md996) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl_996_zhcn, "", "", "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_089, "", "", $md089) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl_089, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrl_089_zhcn, "", "", "089084b7-a469-4f48-856b-98300d6a7fc4.8ddaf840d9f6259ca43465a674b6b03e27c7b072.zh-cn.xlf") | Out-Null

# ==================== de-de sheet ====================
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Cells.Item(2,1).Value = $md996
$ws.Cells.Item(2,2).Value = ".md"
$ws.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,4).Value = "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.de-de.xlf"
$ws.Cells.Item(2,5).Value = "2016-03-18 16:30:55"
$ws.Cells.Item(2,6).Value = $md996
$ws.Cells.Item(2,7).Value = "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.de-de.xlf"
$ws.Cells.Item(2,8).Value = "2016-03-18 16:31:44"
$ws.Cells.Item(2,9).Value = "Include"

$ws.Cells.Item(3,1).Value = $md089
$ws.Cells.Item(3,2).Value = ".md"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "089084b7-a469-4f48-856b-98300d6a7fc4.8ddaf840d9f6259ca43465a674b6b03e27c7b072.de-de.xlf"
$ws.Cells.Item(3,5).Value = "2016-03-18 16:30:26"
$ws.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(3,9).Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_996, "", "", $md996) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl_996, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrl_996_dede, "", "", "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl_996, "", "", $md996) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl_996_dede, "", "", "996f345d-76ed-4a36-aca8-7897bfd7bef4.94c73595a88e2263af8b912a54216edc1884bc38.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_089, "", "", $md089) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl_089, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrl_089_dede, "", "", "089084b7-a469-4f48-856b-98300d6a7fc4.8ddaf840d9f6259ca43465a674b6b03e27c7b072.de-de.xlf") | Out-Null

Write-Host "Done"
